# The underlying OOXML diff for this revision is a pure packaging/
# serialization change: every changed line carries exactly the same set of
# attributes (namespace declarations on <w:document>, and the w:* attributes
# on <w:pgSz>, <w:pgMar>, <w:rFonts>, <w:lang>, <w:latentStyles>,
# <w:lsdException>, <w:style>, <w:tblInd> and <w:tblCellMar>) with no value
# added, removed or changed - only their serialization order differs
# (matching the commit message: "Fixed POI packaging and upgraded to POI
# 3.15", i.e. a newer Apache POI/XMLBeans writer that alphabetizes
# attributes when the fixture was regenerated). There is no document text,
# formatting or structural content to change, so no Word object-model edit
# is required - simply touch the document so the save pipeline runs.
$d = $word.ActiveDocument
$null = $d.Content
